# Refresh the simulation-output figures in the case-study results sheet.
# (Commit: "A bunch of results for the case study are updated")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = [double]"1486.0285581237199"
$ws.Range("H2").Value = [double]"15323.569889185101"
$ws.Range("I2").Value = [double]"23463.267246389201"
$ws.Range("J2").Value = [double]"19364.211036683701"
$ws.Range("K2").Value = [double]"18456.275788321698"
$ws.Range("T2").Value = [double]"0.85199999999999998"
$ws.Range("U2").Value = [double]"0.78200000000000003"
$ws.Range("AA2").Value = [double]"22861.977817288"
$ws.Range("AB2").Value = [double]"21922.1314580617"
$ws.Range("AC2").Value = [double]"26392.876542620099"
$ws.Range("AD2").Value = [double]"22727.9472261546"
$ws.Range("AE2").Value = [double]"23601.375688390701"
$ws.Range("AM2").Value = [double]"0"
$ws.Range("AM2").ClearFormats()
$ws.Range("AN2").Value = [double]"5.0391262418280602E-10"
$ws.Range("AN2").NumberFormat = "0.00E+00"
$ws.Range("AO2").Value = [double]"62.940286557696602"
$ws.Range("AR2").Value = [double]"78093.352518704996"
$ws.Range("AS2").Value = [double]"50684.166873615097"
$ws.Range("AT2").Value = [double]"184.56286151506001"

# Row 3
$ws.Range("G3").Value = [double]"14632.6487718975"
$ws.Range("H3").Value = [double]"27436.348027577002"
$ws.Range("I3").Value = [double]"23783.9556423457"
$ws.Range("J3").Value = [double]"19251.780761506499"
$ws.Range("S3").Value = [double]"0.91700000000000004"
$ws.Range("T3").Value = [double]"0.80800000000000005"
$ws.Range("AA3").Value = [double]"23116.348770770201"
$ws.Range("AB3").Value = [double]"27436.348027577002"
$ws.Range("AC3").Value = [double]"25936.7019000498"
$ws.Range("AD3").Value = [double]"23826.461338498098"
$ws.Range("AL3").Value = [double]"0"
$ws.Range("AM3").Value = [double]"0"
$ws.Range("AN3").Value = [double]"23.992015566669501"
$ws.Range("AR3").Value = [double]"85104.733203326701"
$ws.Range("AS3").Value = [double]"57788.213950614001"
$ws.Range("AT3").Value = [double]"277.22925389091398"

# Row 4
$ws.Range("F4").Value = [double]"13547.725702604401"
$ws.Range("G4").Value = [double]"27415.001906790701"
$ws.Range("H4").Value = [double]"28082.362802842101"
$ws.Range("I4").Value = [double]"22598.365823715001"
$ws.Range("S4").Value = [double]"0.81299999999999994"
$ws.Range("Z4").Value = [double]"27704.9605370233"
$ws.Range("AA4").Value = [double]"27415.001906790701"
$ws.Range("AB4").Value = [double]"28195.1433763475"
$ws.Range("AC4").Value = [double]"27796.267925848799"
$ws.Range("AL4").Value = [double]"3.1445825894581802E-11"
$ws.Range("AM4").Value = [double]"53.423823446857298"
$ws.Range("AR4").Value = [double]"91643.456235952006"
$ws.Range("AS4").Value = [double]"64387.991262001298"
$ws.Range("AT4").Value = [double]"338.283532652492"

# Row 5
$ws.Range("E5").Value = [double]"13547.7252403443"
$ws.Range("F5").Value = [double]"25434.236404408501"
$ws.Range("G5").Value = [double]"32289.541682388801"
$ws.Range("H5").Value = [double]"24438.780256251699"
$ws.Range("Y5").Value = [double]"27704.9595917062"
$ws.Range("Z5").Value = [double]"25434.236404408501"
$ws.Range("AA5").Value = [double]"32289.541682388801"
$ws.Range("AB5").Value = [double]"24937.530873726198"
$ws.Range("AJ5").Value = [double]"0"
$ws.Range("AK5").Value = [double]"0"
$ws.Range("AK5").ClearFormats()
$ws.Range("AL5").Value = [double]"8.8872581741550501E-3"
$ws.Range("AM5").Value = [double]"-5.2927121687274099E-13"
$ws.Range("AN5").Value = [double]"-2.3403928402423801E-11"
$ws.Range("AN5").NumberFormat = "0.00E+00"
$ws.Range("AO5").Value = [double]"2.1306937627062301E-11"
$ws.Range("AO5").NumberFormat = "0.00E+00"
$ws.Range("AR5").Value = [double]"95710.283583393204"
$ws.Range("AS5").Value = [double]"68426.567323889001"
$ws.Range("AT5").Value = [double]"310.03224709915003"

# Row 6
$ws.Range("H6").Value = [double]"68"
$ws.Range("I6").Value = [double]"1126.5994080585399"
$ws.Range("J6").Value = [double]"10093.878331890201"
$ws.Range("K6").Value = [double]"26199.012938838699"
$ws.Range("L6").Value = [double]"19405.187296325199"
$ws.Range("M6").Value = [double]"15934.3495844096"
$ws.Range("T6").Value = [double]"2E-3"
$ws.Range("V6").Value = [double]"0.95"
$ws.Range("W6").Value = [double]"0.96799999999999997"
$ws.Range("X6").Value = [double]"0.78200000000000003"
$ws.Range("Y6").Value = [double]"0.623"
$ws.Range("AG6").Value = [double]"6474.7092417157701"
$ws.Range("AH6").Value = [double]"10625.1350862002"
$ws.Range("AI6").Value = [double]"27065.096011197002"
$ws.Range("AJ6").Value = [double]"24814.8175144823"
$ws.Range("AK6").Value = [double]"25576.805111411901"
$ws.Range("AU6").Value = [double]"0"
$ws.Range("AV6").Value = [double]"14.833988702038001"
$ws.Range("AW6").Value = [double]"79.858087619757001"
$ws.Range("AZ6").Value = [double]"72827.027559522205"
$ws.Range("BA6").Value = [double]"49224.2746097989"
$ws.Range("BB6").Value = [double]"255.731552225672"

# Row 7
$ws.Range("H7").Value = [double]"426.55068922957099"
$ws.Range("I7").Value = [double]"6132.5855215791198"
$ws.Range("J7").Value = [double]"30492.7969034815"
$ws.Range("K7").Value = [double]"25191.975294191299"
$ws.Range("L7").Value = [double]"18499.746147696002"
$ws.Range("W7").Value = [double]"0.93200000000000005"
$ws.Range("X7").Value = [double]"0.80500000000000005"
$ws.Range("AF7").Value = [double]"21327.534461478499"
$ws.Range("AG7").Value = [double]"6151.0386374915997"
$ws.Range("AH7").Value = [double]"30553.904712907301"
$ws.Range("AI7").Value = [double]"27030.016410076401"
$ws.Range("AJ7").Value = [double]"22981.051113601399"
$ws.Range("AU7").Value = [double]"5.7919039222264902E-14"
$ws.Range("AV7").Value = [double]"122.515584251492"
$ws.Range("AW7").Value = [double]"0"
$ws.Range("AW7").ClearFormats()
$ws.Range("AZ7").Value = [double]"80743.654556177498"
$ws.Range("BA7").Value = [double]"57316.646779459297"
$ws.Range("BB7").Value = [double]"431.47672523089102"

# Row 8
$ws.Range("H8").Value = [double]"4971.2507153037004"
$ws.Range("I8").Value = [double]"33312.988828748297"
$ws.Range("J8").Value = [double]"28685.5982541915"
$ws.Range("K8").Value = [double]"22151.230895017899"
$ws.Range("W8").Value = [double]"0.79500000000000004"
$ws.Range("AF8").Value = [double]"4971.2507153037004"
$ws.Range("AG8").Value = [double]"33312.988828748297"
$ws.Range("AH8").Value = [double]"28887.812944805199"
$ws.Range("AI8").Value = [double]"27863.183515745699"
$ws.Range("AU8").Value = [double]"65.038716205687905"
$ws.Range("AZ8").Value = [double]"89121.068693259993"
$ws.Range("BA8").Value = [double]"65572.803108158798"
$ws.Range("BB8").Value = [double]"310.21891684662103"

# Row 9
$ws.Range("G9").Value = [double]"4378.82726051391"
$ws.Range("H9").Value = [double]"33999.999999999898"
$ws.Range("I9").Value = [double]"30701.187446170199"
$ws.Range("J9").Value = [double]"24481.368149927399"
$ws.Range("AE9").Value = [double]"4378.82726051391"
$ws.Range("AF9").Value = [double]"33999.999999999898"
$ws.Range("AG9").Value = [double]"30701.187446170199"
$ws.Range("AH9").Value = [double]"26409.242880180602"
$ws.Range("AT9").Value = [double]"8.67730705067515E-11"
$ws.Range("AZ9").Value = [double]"93561.3828566109"
$ws.Range("BA9").Value = [double]"69976.0807006697"
$ws.Range("BB9").Value = [double]"273.18234600724702"

# Restore the active-cell selection recorded in the saved view
$ws.Range("D15").Select()
